# Issue #44360 - EPPlus can not calculate with Ranges - single Range in Formula
# Adds two new defined names (AGAIN, IF_AGAIN) and formulas in column F / M
# that reference them directly (a bare range-valued defined name used as a
# formula), exercising "calculate with ranges" support.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New defined names used by the new formulas below.
$wb.Names.Add("AGAIN", "=Tabelle1!`$A`$2:`$A`$5")
$wb.Names.Add("IF_AGAIN", "=Tabelle1!`$H`$2:`$H`$6")

# Column F: a single range-valued defined name used directly as a formula.
$ws.Range("F2").Formula = "=AGAIN"
$ws.Range("F3").Formula = "=AGAIN"
$ws.Range("F4").Formula = "=AGAIN"
$ws.Range("F5").Formula = "=AGAIN"

# Column M: same idea, but referencing the other new defined name,
# entered as one block so Excel keeps it as a shared formula.
$ws.Range("M2:M6").Formula = "=IF_AGAIN"

# Move the active selection to F2 (was F14).
$ws.Range("F2").Select()
